$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "CuadroTexto 1") {
        $sh = $cand
    }
}

# --- Resize / reposition the "CuadroTexto 1" textbox ---
$sh.Left   = 33.258160236220476
$sh.Top    = 480.0124122047244
$sh.Width  = 635.0312311023622
$sh.Height = 58.16256968503937

# --- Update the text of the first paragraph ---
$tf = $sh.TextFrame
$tr = $tf.TextRange
$chars = $tr.Characters(27, 4)
$r3 = $chars.InsertAfter(")")
$r2 = $chars.InsertAfter("etc")
$r1 = $chars.InsertAfter(" (Se requieren materiales como leds ")

# --- Add a new empty paragraph at the end ---
$len = $tr.Text.Length
$endRange = $tr.Characters($len, 0)
$newPara = $endRange.InsertAfter("`r")
